# Applies the diff to poster/292C Project Poster.pptx:
#  - Reposition the "Tool Architecture" banner shape (id 41)
#  - Reposition the "Goal: given a puzzle..." body shape (id 60)
#  - Reposition/resize the architecture diagram picture (id 13)
#  - Merge two text runs in the pSUS description shape (id 58)

function Get-ShapeById {
    param($slide, $id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point conversion (1 pt = 12700 EMU), as COM Shape Left/Top/Width/Height are in points.
$emuPerPt = 12700.0

# --- Shape id 41: "Tool Architecture" header bar ---
$shape1 = Get-ShapeById $s 41
$shape1.Left = 503611 / $emuPerPt
$shape1.Top = 14795288 / $emuPerPt

# --- Shape id 60: "Goal: given a puzzle with holes..." body text box ---
$shape2 = Get-ShapeById $s 60
$shape2.Left = 134144 / $emuPerPt
$shape2.Top = 10652759 / $emuPerPt

# --- Shape id 13: architecture diagram picture ---
$shape3 = Get-ShapeById $s 13
$shape3.Left = 1148894 / $emuPerPt
$shape3.Top = 15676682 / $emuPerPt
$shape3.Width = 12333908 / $emuPerPt
$shape3.Height = 14785092 / $emuPerPt

# --- Shape id 58: merge " compiles puzzle to " and "SMT formula, passes to Z3" runs ---
$shape4 = Get-ShapeById $s 58
$tr = $shape4.TextFrame.TextRange
$full = $tr.Text
$tailLen = $full.Length - 4
$tail = $tr.Characters(5, $tailLen)
$tail.Text = " compiles puzzle to SMT formula, passes to Z3"
